$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts existing rows 35-105 down to 36-106
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record's data.
# Columns A, B, C, E, F, G, H, N, O, Q, R mirror the neighboring records (unchanged template),
# while D, I, J, K, L, M, P hold the new values for this record.
$ws.Cells.Item(35, 1).Value = 1
$ws.Cells.Item(35, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(35, 4).Value = 44967
$ws.Cells.Item(35, 5).Value = 15
$ws.Cells.Item(35, 6).Value = 100112040
$ws.Cells.Item(35, 7).Value = "Cilantro"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Segunda"
$ws.Cells.Item(35, 10).Value = 450
$ws.Cells.Item(35, 11).Value = 3000
$ws.Cells.Item(35, 12).Value = 3500
$ws.Cells.Item(35, 13).Value = 3222
$ws.Cells.Item(35, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(35, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value = 1611
$ws.Cells.Item(35, 17).Value = 2
$ws.Cells.Item(35, 18).Value = "Hortaliza"
